$d = $word.ActiveDocument

# --- 1) First numbered paragraph (table 1, "Практична підготовка") ---
# Drop the pStyle/numPr auto-numbering and splice in a literal "15. " run
# in front of the existing text run.
$p1 = $d.Tables.Item(1).Cell(1, 1).Range.Paragraphs.Item(1)
$xml1 = '<w:p w14:paraId="76B7BC5E" w14:textId="77777777" w:rsidR="003E3F51" w:rsidRPr="00DD2FD6" w:rsidRDefault="003E3F51" w:rsidP="005270E2"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="uk-UA"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve">15. </w:t></w:r><w:r w:rsidRPr="00DD2FD6"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="uk-UA"/></w:rPr><w:t>Практична підготовка</w:t></w:r></w:p>'
$p1.Range.InsertXML($xml1)

# --- 2) Drop the stray "_GoBack" bookmark next to the "#P" placeholder ---
$t3 = $d.Tables.Item(3)
$pb = $t3.Cell(3, 3).Range.Paragraphs.Item(1)
$xml2 = '<w:p w14:paraId="00AD4A70" w14:textId="3B611D86" w:rsidR="00CD5D9C" w:rsidRPr="00DD2FD6" w:rsidRDefault="00E1283A" w:rsidP="003E3F51"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00633508" w:rsidRPr="00DD2FD6"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr><w:t>#P</w:t></w:r></w:p>'
$pb.Range.InsertXML($xml2)

# --- 3) Second numbered paragraph (table 4, "Атестація") ---
# Drop the pStyle/numPr auto-numbering, splice in literal "16" + the
# relocated "_GoBack" bookmark + ". " ahead of the existing text run.
$p75 = $d.Tables.Item(4).Cell(1, 1).Range.Paragraphs.Item(1)
$xml3 = '<w:p w14:paraId="159070A0" w14:textId="77777777" w:rsidR="008F1B73" w:rsidRPr="00DD2FD6" w:rsidRDefault="008F1B73" w:rsidP="008F1B73"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="uk-UA"/></w:rPr><w:t>16</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="uk-UA"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidRPr="00DD2FD6"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="uk-UA"/></w:rPr><w:t>Атестація</w:t></w:r></w:p>'
$p75.Range.InsertXML($xml3)

Write-Host "edits applied"
